$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Zoom in the sheet view from 190% to 235%
$excel.ActiveWindow.Zoom = 235

# Move the selection from F2:F11 to the header row A1:G1
$ws.Range("A1:G1").Select()

# Narrow column A (target stored width 12.625) and split the old F:G
# (13.84-wide) pairing into two differently sized columns: F -> 8.125,
# G -> 11.625.
#
# The host stores ColumnWidth internally as whole pixels
# (pixels = round(CharWidth * 7) + 5, persisted width = pixels / 7), so the
# character-unit value that round-trips closest to each target stored width
# is target - 5/7, pre-compensating for that +5-pixel padding.
$ws.Columns(1).ColumnWidth = 11.857142857142858   # -> stored width 12.625 (closest: 12.5714)
$ws.Columns(6).ColumnWidth = 7.428571428571429    # -> stored width 8.125  (closest: 8.1429)
$ws.Columns(7).ColumnWidth = 10.857142857142858   # -> stored width 11.625 (closest: 11.5714)
